$wb = $excel.ActiveWorkbook

$values = @{
  "run_1" = @{
    2 = 30.25177407264709
    3 = 30.05023574829102
    4 = 30.02302861213684
    5 = 29.9938018321991
    6 = 29.97183775901794
    7 = 29.96538615226746
    8 = 29.88832664489746
    9 = 29.82996368408203
    10 = 29.97231841087341
    11 = 30.4281747341156
    12 = 29.95809626579285
    13 = 29.89967465400696
    14 = 29.93186545372009
    15 = 29.99841547012329
    16 = 29.8017270565033
    17 = 29.88027691841125
    18 = 29.83977055549622
    19 = 29.9912941455841
    20 = 29.90792393684387
    21 = 30.32831263542176
  }
  "run_2" = @{
    2 = 30.17731547355652
    3 = 30.03387403488159
    4 = 29.96002197265625
    5 = 30.05636024475098
    6 = 30.01436924934387
    7 = 30.10537576675415
    8 = 29.94202184677124
    9 = 29.91912245750427
    10 = 30.04311943054199
    11 = 30.3587737083435
    12 = 30.05744576454162
    13 = 30.09493708610535
    14 = 30.0193338394165
    15 = 30.10676121711731
    16 = 30.0355212688446
    17 = 30.0868775844574
    18 = 29.95192146301269
    19 = 30.01825380325317
    20 = 29.88979983329773
    21 = 30.312096118927
  }
  "run_3" = @{
    2 = 30.25078916549682
    3 = 29.97223806381226
    4 = 29.89531278610229
    5 = 30.13256931304932
    6 = 29.98906469345093
    7 = 30.02297377586365
    8 = 29.85321617126465
    9 = 29.82651686668396
    10 = 29.74919700622558
    11 = 30.25184345245361
    12 = 29.9210352897644
    13 = 29.87582468986511
    14 = 29.96723008155823
    15 = 29.84311985969543
    16 = 29.82358169555664
    17 = 30.0328266620636
    18 = 29.91261529922485
    19 = 29.91936635971069
    20 = 29.98606085777283
    21 = 30.21571969985962
  }
  "run_4" = @{
    2 = 30.20046782493592
    3 = 30.01102018356323
    4 = 29.79540157318115
    5 = 29.95803785324097
    6 = 29.86515045166016
    7 = 30.03622674942017
    8 = 29.92337989807129
    9 = 29.96232914924622
    10 = 29.83899784088135
    11 = 30.20390391349792
    12 = 29.92948460578918
    13 = 29.94855785369873
    14 = 29.92574310302734
    15 = 30.07432675361633
    16 = 29.95177435874939
    17 = 29.94523763656616
    18 = 29.92781734466553
    19 = 30.04430389404297
    20 = 30.14005446434021
    21 = 30.15230369567871
  }
  "run_5" = @{
    2 = 30.2267017364502
    3 = 29.89012861251831
    4 = 29.82972598075867
    5 = 29.91176557540894
    6 = 29.90625834465027
    7 = 30.12267518043518
    8 = 29.76196551322937
    9 = 29.77659964561462
    10 = 30.12089443206787
    11 = 30.16909146308899
    12 = 29.92731475830078
    13 = 29.82842421531677
    14 = 29.86576271057129
    15 = 29.85349011421204
    16 = 29.86205983161926
    17 = 29.94987511634827
    18 = 29.88572072982788
    19 = 29.92448139190674
    20 = 30.06937885284424
    21 = 30.11360836029053
  }
}

foreach ($sheetName in $values.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)
  $rowValues = $values[$sheetName]
  foreach ($row in $rowValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $rowValues[$row]
  }
}
